$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1543.2222
$ws.Range("J17").Value = 1543.2222
$ws.Range("L17").Value = 4629.6666
$ws.Range("N17").Value = -4965.6666
$ws.Range("H19").Value = 2087.0908
$ws.Range("J19").Value = 2807.25
$ws.Range("L19").Value = 2807.25
$ws.Range("N19").Value = -3157.25
$ws.Range("H74").Value = 2332.3333
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 2332.3333
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null
$ws.Range("H113").Value = 38464860
$ws.Range("I113").Value = 12502715
$ws.Range("K113").Value = 12502715
$ws.Range("M113").Value = -12499461
$ws.Range("H135").Value = 2230.28
$ws.Range("I135").Value = 2184.9546
$ws.Range("K135").Value = 19664.5914
$ws.Range("M135").Value = -17129.5914
$ws.Range("H137").Value = 11225.692
$ws.Range("I137").Value = 6363.1113
$ws.Range("K137").Value = 19089.3339
$ws.Range("M137").Value = -16539.3339

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 31252544
$ws.Range("I45").Value = 55557120
$ws.Range("K45").Value = 55557120
$ws.Range("M45").Value = -55556743
$ws.Range("H61").Value = 24201356
$ws.Range("I61").Value = 18525614
$ws.Range("K61").Value = 18525614
$ws.Range("M61").Value = -18525402
$ws.Range("H132").Value = 4071.9167
$ws.Range("I132").Value = 2116.96
$ws.Range("J132").Value = 8515
$ws.Range("K132").Value = 6350.88
$ws.Range("L132").Value = 25545
$ws.Range("M132").Value = -3820.88
$ws.Range("N132").Value = -30605
$ws.Range("H136").Value = 24201356
$ws.Range("I136").Value = 18525614
$ws.Range("K136").Value = 55576842
$ws.Range("M136").Value = -55574292

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1137.4242
$ws.Range("I94").Value = 1085.409
$ws.Range("J94").Value = 1241.4546
$ws.Range("K94").Value = 1085.409
$ws.Range("L94").Value = 1241.4546
$ws.Range("M94").Value = -634.4090000000001
$ws.Range("N94").Value = -2143.4546
$ws.Range("H134").Value = 2501409.5
$ws.Range("I134").Value = 1874.6666
$ws.Range("J134").Value = 10000014
$ws.Range("K134").Value = 5623.9998
$ws.Range("L134").Value = 30000042
$ws.Range("M134").Value = -3088.9998
$ws.Range("N134").Value = -30005112

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -888
$ws.Range("H23").Value = 1750
$ws.Range("I23").Value = 500
$ws.Range("K23").Value = 500
$ws.Range("M23").Value = -260
$ws.Range("H27").Value = 1750
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 500
$ws.Range("M27").Value = -308
$ws.Range("H50").Value = 45197.8
$ws.Range("J50").Value = 45197.8
$ws.Range("L50").Value = 45197.8
$ws.Range("N50").Value = -46447.8
$ws.Range("H132").Value = 1968.3846
$ws.Range("I132").Value = 1968.3846
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5905.1538
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3375.1538
$ws.Range("N132").Value = $null

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4139.6875
$ws.Range("I80").Value = 3402
$ws.Range("J80").Value = 4188.8667
$ws.Range("K80").Value = 10206
$ws.Range("L80").Value = 12566.6001
$ws.Range("M80").Value = -9270
$ws.Range("N80").Value = -14438.6001
$ws.Range("H83").Value = 4139.6875
$ws.Range("I83").Value = 3402
$ws.Range("J83").Value = 4188.8667
$ws.Range("K83").Value = 30618
$ws.Range("L83").Value = 37699.8003
$ws.Range("M83").Value = -25938
$ws.Range("N83").Value = -47059.8003
$ws.Range("H107").Value = 603.8788
$ws.Range("J107").Value = 856.53845
$ws.Range("L107").Value = 2569.61535
$ws.Range("N107").Value = -6409.61535
$ws.Range("H113").Value = 1401.2222
$ws.Range("J113").Value = 1817
$ws.Range("L113").Value = 5451
$ws.Range("N113").Value = -9791
$ws.Range("H115").Value = 70031
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 70031
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 210093
$ws.Range("M115").Value = $null
$ws.Range("N115").Value = -212443
$ws.Range("I131").Value = 3541.8823
$ws.Range("J131").Value = 4967.1704
$ws.Range("K131").Value = 10625.6469
$ws.Range("L131").Value = 14901.5112
$ws.Range("M131").Value = -5585.6469
$ws.Range("N131").Value = -24981.5112
$ws.Range("H132").Value = 1655.2
$ws.Range("I132").Value = 902.6667
$ws.Range("K132").Value = 8124.0003
$ws.Range("M132").Value = -5594.0003
$ws.Range("H133").Value = 6706
$ws.Range("I133").Value = 6265
$ws.Range("K133").Value = 18795
$ws.Range("M133").Value = -13735

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 120000
$ws.Range("J62").Value = 120000
$ws.Range("L62").Value = 120000
$ws.Range("N62").Value = -121372
$ws.Range("H63").Value = 26103.5
$ws.Range("J63").Value = 26103.5
$ws.Range("L63").Value = 26103.5
$ws.Range("N63").Value = -27475.5
$ws.Range("H65").Value = 120000
$ws.Range("J65").Value = 120000
$ws.Range("L65").Value = 360000
$ws.Range("N65").Value = -366864
$ws.Range("H66").Value = 26103.5
$ws.Range("J66").Value = 26103.5
$ws.Range("L66").Value = 78310.5
$ws.Range("N66").Value = -85174.5
$ws.Range("H70").Value = 7577.4
$ws.Range("I70").Value = 6995.6665
$ws.Range("J70").Value = 8450
$ws.Range("K70").Value = 6995.6665
$ws.Range("L70").Value = 8450
$ws.Range("M70").Value = -6725.6665
$ws.Range("N70").Value = -8990
$ws.Range("H73").Value = 7577.4
$ws.Range("I73").Value = 6995.6665
$ws.Range("J73").Value = 8450
$ws.Range("K73").Value = 6995.6665
$ws.Range("L73").Value = 8450
$ws.Range("M73").Value = -6059.6665
$ws.Range("N73").Value = -10322
$ws.Range("H109").Value = 45216.75
$ws.Range("J109").Value = 45216.75
$ws.Range("L109").Value = 45216.75
$ws.Range("N109").Value = -47296.75
$ws.Range("H122").Value = 2820
$ws.Range("I122").Value = 2775
$ws.Range("K122").Value = 8325
$ws.Range("M122").Value = -5875
$ws.Range("H132").Value = 33340974
$ws.Range("I132").Value = 43482720
$ws.Range("J132").Value = 18093.572
$ws.Range("K132").Value = 130448160
$ws.Range("L132").Value = 54280.716
$ws.Range("M132").Value = -130445630
$ws.Range("N132").Value = -59340.716

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 73601.60000000001
$ws.Range("I7").Value = 5875.625
$ws.Range("J7").Value = 151002.72
$ws.Range("K7").Value = 5875.625
$ws.Range("L7").Value = 151002.72
$ws.Range("M7").Value = -5763.625
$ws.Range("N7").Value = -151226.72
$ws.Range("H126").Value = 73601.60000000001
$ws.Range("I126").Value = 5875.625
$ws.Range("J126").Value = 151002.72
$ws.Range("K126").Value = 17626.875
$ws.Range("L126").Value = 453008.16
$ws.Range("M126").Value = -15156.875
$ws.Range("N126").Value = -457948.16
$ws.Range("H132").Value = 426522.75
$ws.Range("I132").Value = 11167.15
$ws.Range("K132").Value = 33501.45
$ws.Range("M132").Value = -30971.45

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 26646.334
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 26646.334
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 26646.334
$ws.Range("M39").Value = $null
$ws.Range("N39").Value = -27472.334
$ws.Range("H42").Value = 58888
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 58888
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 58888
$ws.Range("M42").Value = $null
$ws.Range("N42").Value = -59644
$ws.Range("H122").Value = 1383.7222
$ws.Range("I122").Value = 1369.375
$ws.Range("K122").Value = 4108.125
$ws.Range("M122").Value = -1658.125
